$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cycle counts in the raw data table (col A) ---
$ws.Range("A3").Value = 143164482
$ws.Range("A4").Value = 143013479
$ws.Range("A5").Value = 121658327
$ws.Range("A7").Value = 116110440
$ws.Range("A8").Value = 63360132
$ws.Range("A9").Value = 50182384
$ws.Range("A14").Value = 5633800

# --- Fix a method name (dropped the $CmpComplex inner-class qualifier) ---
$ws.Range("C7").Value = "CoreListJoinA.calc_func"

# --- Add a note explaining why some rows moved out of the main table ---
$ws.Range("M5").Value = "Before switching to coremk_ch2 as main version, these were not inlined"

# --- Update the "last updated" stamp ---
$ws.Range("J3").Value = "UDPATED 20180327"

# --- Move the CmpComplex.compare and CmpIdx.compare rows out of the main table ---
$ws.Range("A6:C6").Cut($ws.Range("M6"))
$ws.Range("A16:C16").Cut($ws.Range("M16"))

# --- Recolor the highlight fill used by J3:K3 from orange to red ---
$ws.Range("J3:K3").Interior.Color = 255

# --- Restore the previous selection ---
$ws.Range("L27").Select()
